$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 235; existing rows 235:313 shift down to 236:314.
$ws.Rows("235:235").Insert()

# Populate the newly inserted row 235 with this week's data point.
$ws.Cells.Item(235, 1).Value = 10
$ws.Cells.Item(235, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(235, 3).Value = "La Araucanía"
$ws.Cells.Item(235, 4).Value = 45120
$ws.Cells.Item(235, 5).Value = 9
$ws.Cells.Item(235, 6).Value = 100112005
$ws.Cells.Item(235, 7).Value = "Puerro"
$ws.Cells.Item(235, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 100
$ws.Cells.Item(235, 11).Value = 9000
$ws.Cells.Item(235, 12).Value = 9000
$ws.Cells.Item(235, 13).Value = 9000
$ws.Cells.Item(235, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(235, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(235, 16).Value = 750
$ws.Cells.Item(235, 17).Value = 12
$ws.Cells.Item(235, 18).Value = "Hortaliza"
